# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data using K instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 2
